# Fixing errors in example upload files.

$wb = $excel.ActiveWorkbook

# --- Service Contacts sheet: widen column A, move selection ---
$wsSvc = $wb.Worksheets.Item("Service Contacts")
$wsSvc.Columns.Item(1).ColumnWidth = 13.666666666666666
$wsSvc.Activate()
[void]$wsSvc.Range("D3").Select()

# --- Practitioners sheet: widen columns, add new data row, move selection ---
$wsPrac = $wb.Worksheets.Item("Practitioners")
$wsPrac.Activate()

$wsPrac.Columns.Item(1).ColumnWidth = 13.833333333333334
$wsPrac.Columns.Item(3).ColumnWidth = 12.166666666666666
$wsPrac.Columns.Item(6).ColumnWidth = 12

$wsPrac.Range("A6").Value = "PHN999:NFP02"
$wsPrac.Range("B6").Value = "P01"
$wsPrac.Range("C6").Value = 8
$wsPrac.Range("D6").Value = 1
$wsPrac.Range("E6").Value = 1973
$wsPrac.Range("F6").Value = 2
$wsPrac.Range("G6").Value = 1
$wsPrac.Range("H6").Value = 1
$wsPrac.Range("I6").Value = "tag1"

[void]$wsPrac.Columns.Item(7).Select()

$wb.Save()
